$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds crypto prices and column E holds 1h % change, both stored
# as plain text (inlineStr) in the sheet. Apply the refreshed figures.

$ws.Range("D2").Value = "64.912.11"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").Value = "2.651.61"
$ws.Range("E3").Value = "  +2.90%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.66%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("E9").Value = "  +7.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.399"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000185"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +19.32%  "
$ws.Range("D15").Value = "3.130.24"
$ws.Range("D16").Value = "64.862.87"
$ws.Range("E16").Value = "  +2.82%  "
$ws.Range("D17").Value = "2.647.48"
$ws.Range("E17").Value = "  +3.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.78%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.164"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "0.0₃0942"
$ws.Range("E30").Value = "  +10.16%  "
$ws.Range("E31").Value = "  +4.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "520.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.426"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "165.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0614"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.648"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.45%  "
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0984"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.67%  "
